# Update benchmark: 2026-01-29 06:56:55 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "39,87 TRY - 79,76 TRY - 797,68 TRY"
$ws.Range("G4").Value = "27,84 TRY - 55,69 TRY - 398,83 TRY"
$ws.Range("G5").Value = "7,97 TRY - 15,96 TRY - 199,41 TRY"

$ws.Range("D6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("G6").Value = "8.300,01 TL - 99,71 TL"

$ws.Range("G8").Value = "19,94 TRY - 39,88 TRY - 398,84 TRY"
$ws.Range("G9").Value = "13,92 TRY - 27,85 TRY - 199,42 TRY"
$ws.Range("G10").Value = "3,99 TRY - 7,98 TRY - 99,71 TRY"
$ws.Range("G11").Value = "3,99 TRY - 7,98 TRY - 99,71 TRY"

$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"

$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"

$ws.Range("D14").Value = "3.500 TL - 13.500 TL"
$ws.Range("G14").Value = "8.300 TL - 7,97 TL"
